# Atualizacao de bases das ligas, do dia: 09-03-2024 as 13:07
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 73 and 74: the two fixtures were swapped (everything except the
# running index in column A stays put; all other columns exchange values). ---

$ws.Range("B73").Value = 7646749
$ws.Range("F73").Value = "Brisbane Roar"
$ws.Range("G73").Value = "Newcastle Jets"
$ws.Range("I73").Value = 2
$ws.Range("J73").Value = "H"
$ws.Range("K73").Value = 1.909
$ws.Range("L73").Value = 4
$ws.Range("M73").Value = 3.4
$ws.Range("N73").Value = 2.4
$ws.Range("O73").Value = 4
$ws.Range("P73").Value = 2.6
$ws.Range("Q73").Value = 0
$ws.Range("R73").Value = 1.83
$ws.Range("S73").Value = 2.07
$ws.Range("T73").Value = 3.25
$ws.Range("U73").Value = 1.9
$ws.Range("V73").Value = 1.95
$ws.Range("W73").Value = 1.4
$ws.Range("Y73").Value = -1
$ws.Range("Z73").Value = 0.8300000000000001
$ws.Range("AA73").Value = -1
$ws.Range("AB73").Value = 0.8999999999999999

$ws.Range("B74").Value = 7646750
$ws.Range("F74").Value = "Perth Glory"
$ws.Range("G74").Value = "Wellington Phoenix"
$ws.Range("I74").Value = 4
$ws.Range("J74").Value = "A"
$ws.Range("K74").Value = 2.45
$ws.Range("L74").Value = 3.75
$ws.Range("M74").Value = 2.55
$ws.Range("N74").Value = 3.1
$ws.Range("O74").Value = 3.8
$ws.Range("P74").Value = 2.05
$ws.Range("Q74").Value = 0.25
$ws.Range("R74").Value = 2
$ws.Range("S74").Value = 1.85
$ws.Range("T74").Value = 3
$ws.Range("U74").Value = 1.925
$ws.Range("V74").Value = 1.925
$ws.Range("W74").Value = -1
$ws.Range("Y74").Value = 1.05
$ws.Range("Z74").Value = -1
$ws.Range("AA74").Value = 0.8500000000000001
$ws.Range("AB74").Value = 0.925

# --- Rows 120-122: matches have now been played, so the final score
# (H/I) and result code (J) are filled in and the closing odds / P&L
# columns are refreshed. ---

$ws.Range("H120").Value = 1
$ws.Range("I120").Value = 3
$ws.Range("J120").Value = "A"
$ws.Range("N120").Value = 1.727
$ws.Range("O120").Value = 4.5
$ws.Range("R120").Value = 1.95
$ws.Range("S120").Value = 1.9
$ws.Range("U120").Value = 1.975
$ws.Range("V120").Value = 1.875
$ws.Range("W120").Value = -1
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = 3
$ws.Range("Z120").Value = -1
$ws.Range("AA120").Value = 0.8999999999999999
$ws.Range("AB120").Value = 0.9750000000000001
$ws.Range("AC120").Value = -1

$ws.Range("H121").Value = 1
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = "H"
$ws.Range("N121").Value = 2
$ws.Range("O121").Value = 3.75
$ws.Range("P121").Value = 3.5
$ws.Range("R121").Value = 2.025
$ws.Range("S121").Value = 1.825
$ws.Range("U121").Value = 1.975
$ws.Range("V121").Value = 1.875
$ws.Range("W121").Value = 1
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = -1
$ws.Range("Z121").Value = 1.025
$ws.Range("AA121").Value = -1
$ws.Range("AB121").Value = -1
$ws.Range("AC121").Value = 0.875

$ws.Range("H122").Value = 1
$ws.Range("I122").Value = 2
$ws.Range("J122").Value = "A"
$ws.Range("N122").Value = 4
$ws.Range("O122").Value = 3.8
$ws.Range("P122").Value = 1.85
$ws.Range("Q122").Value = 0.5
$ws.Range("R122").Value = 2.01
$ws.Range("S122").Value = 1.89
$ws.Range("T122").Value = 3.25
$ws.Range("U122").Value = 2
$ws.Range("V122").Value = 1.85
$ws.Range("W122").Value = -1
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = 0.8500000000000001
$ws.Range("Z122").Value = -1
$ws.Range("AA122").Value = 0.8899999999999999
$ws.Range("AB122").Value = -0.5
$ws.Range("AC122").Value = 0.425

# --- Rows 123-125: still upcoming fixtures, but the id / kickoff-time /
# opening-odds refresh moved each one down to the next scheduled match. ---

$ws.Range("B123").Value = 7128012
$ws.Range("E123").Value = 45361.125
$ws.Range("F123").Value = "Macarthur FC"
$ws.Range("G123").Value = "Central Coast Mariners"
$ws.Range("K123").Value = 2.4
$ws.Range("L123").Value = 3.5
$ws.Range("M123").Value = 2.75
$ws.Range("N123").Value = 3.1
$ws.Range("O123").Value = 3.5
$ws.Range("P123").Value = 2.25
$ws.Range("Q123").Value = 0.25
$ws.Range("R123").Value = 1.97
$ws.Range("S123").Value = 1.93
$ws.Range("T123").Value = 3
$ws.Range("U123").Value = 2
$ws.Range("V123").Value = 1.85

$ws.Range("P124").Value = 4.75
$ws.Range("R124").Value = 1.92
$ws.Range("S124").Value = 1.98

$ws.Range("B125").Value = 7662592
$ws.Range("E125").Value = 45363.20833333334
$ws.Range("F125").Value = "Melbourne City"
$ws.Range("G125").Value = "Western Sydney Wanderers"
$ws.Range("K125").Value = 2.1
$ws.Range("L125").Value = 3.6
$ws.Range("M125").Value = 3.25
$ws.Range("N125").Value = 2.1
$ws.Range("O125").Value = 3.6
$ws.Range("P125").Value = 3.25
$ws.Range("Q125").Value = -0.25
$ws.Range("R125").Value = 1.84
$ws.Range("S125").Value = 2.06
$ws.Range("U125").Value = 1.825
$ws.Range("V125").Value = 2.025

# --- New fixtures appended at rows 126 and 127. Copy the format from the
# last existing data row first so the index (A) and date (E) columns keep
# the same style (bold/bordered index, date-formatted kickoff time). ---

$ws.Range("A125:V125").Copy()
$ws.Range("A126:V126").PasteSpecial(-4122)
$ws.Range("A127:V127").PasteSpecial(-4122)

$ws.Range("A126").Value = 124
$ws.Range("B126").Value = 7127392
$ws.Range("C126").Value = "Australia ALeague"
$ws.Range("D126").Value = "Australia ALeague"
$ws.Range("E126").Value = 45365.20833333334
$ws.Range("F126").Value = "Western United FC"
$ws.Range("G126").Value = "Melbourne Victory"
$ws.Range("K126").Value = 4.5
$ws.Range("L126").Value = 3.5
$ws.Range("M126").Value = 1.8
$ws.Range("N126").Value = 4.2
$ws.Range("O126").Value = 3.5
$ws.Range("P126").Value = 1.909
$ws.Range("Q126").Value = 0.5
$ws.Range("R126").Value = 1.95
$ws.Range("S126").Value = 1.95
$ws.Range("T126").Value = 3
$ws.Range("U126").Value = 2.025
$ws.Range("V126").Value = 1.825
$ws.Range("W126").Value = 0
$ws.Range("X126").Value = 0
$ws.Range("Y126").Value = 0
$ws.Range("Z126").Value = 0
$ws.Range("AA126").Value = 0

$ws.Range("A127").Value = 125
$ws.Range("B127").Value = 7127389
$ws.Range("C127").Value = "Australia ALeague"
$ws.Range("D127").Value = "Australia ALeague"
$ws.Range("E127").Value = 45366.23958333334
$ws.Range("F127").Value = "Newcastle Jets"
$ws.Range("G127").Value = "Adelaide United"
$ws.Range("K127").Value = 2.25
$ws.Range("L127").Value = 3.5
$ws.Range("M127").Value = 3
$ws.Range("N127").Value = 2.45
$ws.Range("O127").Value = 3.5
$ws.Range("P127").Value = 2.75
$ws.Range("Q127").Value = 0
$ws.Range("R127").Value = 1.83
$ws.Range("S127").Value = 2.07
$ws.Range("T127").Value = 3.25
$ws.Range("U127").Value = 1.85
$ws.Range("V127").Value = 2
$ws.Range("W127").Value = 0
$ws.Range("X127").Value = 0
$ws.Range("Y127").Value = 0
$ws.Range("Z127").Value = 0
$ws.Range("AA127").Value = 0
